# Build site at 2022-09-26 16:07:08 UTC
# This workbook is a single-sheet "course syllabus" table laid out as
# Col A = field label, Col B = value, Col C = value (highlighted copy).
# The edit removes the standalone row that held "5701460 - Antonio Iacono"
# (old row 13, which had no label in column A), shifting every row below it
# up by one, and then refreshes a handful of cells whose underlying data
# changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the old row 13 ("5701460 - Antonio Iacono" with no A-column
#    label) entirely - this shifts rows 14-24 up to become rows 13-23,
#    carrying their row heights along with them.
$ws.Rows.Item(13).Delete()

# 2) Now patch the handful of cells whose text content changed as part of
#    this edit (after the shift above).

# Row 10 ("Objetivos:") value now holds the professor info instead of the
# old course-objectives paragraph.
$ws.Range("B10").Value = "5701460 - Antonio Iacono"
$ws.Range("C10").Value = "5701460 - Antonio Iacono"

# Row 13 ("Programa resumido:") value is now "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 ("Programa:") value is now the activation date.
$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

# Row 18 ("Método:") value now holds the professor info.
$ws.Range("B18").Value = "5701460 - Antonio Iacono"
$ws.Range("C18").Value = "5701460 - Antonio Iacono"

# Row 19 ("Critério:") value is now "Provas e Trabalhos".
$ws.Range("B19").Value = "Provas e Trabalhos"
$ws.Range("C19").Value = "Provas e Trabalhos"
